$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove D2 and D3 entirely (historical_growth_revenue_last_5_years)
$ws.Range("D2").ClearContents()
$ws.Range("D3").ClearContents()

foreach ($row in 2, 3) {
    $ws.Range("G$row").Value = -0.1003717472118959
    $ws.Range("H$row").Value = -0.1003717472118959
    $ws.Range("I$row").Value = -0.254275092936803
    $ws.Range("J$row").Value = -0.254275092936803
    $ws.Range("K$row").Value = -0.8129999999999999
    $ws.Range("L$row").Value = -0.3022304832713755

    $ws.Range("U$row").Value = 0.075
    $ws.Range("V$row").Value = 0.009566326530612245
    $ws.Range("W$row").Value = -0.09713261648745521
    $ws.Range("X$row").Value = 0.08318406724216443
    $ws.Range("Y$row").Value = -0.1803166837296196
    $ws.Range("Z$row").Value = 0.2778925619834711
    $ws.Range("AA$row").Value = -0.07066115702479339
    $ws.Range("AB$row").Value = 0.07910046291419674
    $ws.Range("AC$row").Value = -0.1497616199389901
    $ws.Range("AD$row").Value = 1.33
    $ws.Range("AF$row").Value = 1.33
    $ws.Range("AG$row").Value = 1.255
    $ws.Range("AH$row").Value = 0.1450381679389313
    $ws.Range("AI$row").Value = 0.1631901840490798
    $ws.Range("AJ$row").Value = 0.1379879054425509
    $ws.Range("AK$row").Value = 0.1554179566563468
    $ws.Range("AL$row").Value = 0.08699999999999999
    $ws.Range("AM$row").Value = 0.08699999999999999
    $ws.Range("AN$row").Value = -4.130434782608695
    $ws.Range("AO$row").Value = -7.862068965517243
    $ws.Range("AP$row").Value = -3.897515527950311
    $ws.Range("AQ$row").Value = -7.862068965517243
}
